$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("James", "Thomas", "jessica84@example.org", "253.302.7370")
    ,@("William", "Roth", "billy12@example.com", "873-841-3042x60978")
    ,@("Roy", "Lyons", "monica41@example.com", "4085959353")
    ,@("Russell", "Thompson", "akennedy@example.org", "384-701-4754")
    ,@("Christine", "Smith", "bhodges@example.net", "001-868-696-3107x56462")
    ,@("Michael", "Gardner", "sherryprice@example.org", "821-226-3604")
    ,@("Timothy", "Perez", "danielfoster@example.org", "759.785.4635")
    ,@("Thomas", "Martinez", "seansanchez@example.net", "315.704.4982x7779")
    ,@("Robert", "Castillo", "wwalker@example.org", "(961)298-0738x92660")
    ,@("Sharon", "Miller", "vphillips@example.net", "9189362930")
    ,@("Mia", "Liu", "greenemaria@example.com", "(830)928-9081")
    ,@("David", "Downs", "dgordon@example.org", "+1-903-513-8724x6346")
    ,@("Brandy", "White", "kbenson@example.com", "4184204038")
    ,@("Jennifer", "Allen", "ysmith@example.com", "+1-979-482-8665x201")
    ,@("Amy", "Walker", "calebpowell@example.com", "+1-479-445-3207x5997")
    ,@("Jessica", "Good", "adam35@example.com", "001-303-830-1826x617")
    ,@("Cynthia", "Freeman", "kevinsharp@example.org", "+1-639-851-0189x756")
    ,@("Joshua", "May", "kirsten47@example.net", "(848)705-3659x100")
    ,@("Heidi", "Smith", "jimmy04@example.org", "471-669-3904")
    ,@("Susan", "Francis", "michael26@example.net", "335-426-7572")
    ,@("Kathleen", "Barnes", "jordanchelsea@example.com", "+1-976-854-0346x6336")
    ,@("Julie", "Hurst", "marklandry@example.net", "804.494.4234x10343")
    ,@("Linda", "Odonnell", "kevin43@example.org", "329.608.5203x82743")
    ,@("Yolanda", "Keller", "mary86@example.org", "521-769-7755")
    ,@("Sarah", "Alvarado", "toniwright@example.com", "859.879.3812")
    ,@("Ricky", "King", "wilsondavid@example.org", "688-978-3814")
    ,@("Katherine", "Clark", "zmccarthy@example.org", "253-872-0966")
    ,@("Lucas", "King", "whitneyferguson@example.com", "932-783-9197x072")
    ,@("Laura", "Andrews", "thomasoneill@example.org", "+1-783-289-4694x68221")
    ,@("Robert", "Peterson", "susanthomas@example.com", "+1-507-227-3275")
    ,@("Jonathan", "George", "rosecox@example.com", "001-632-639-6738")
    ,@("Gabriela", "Campbell", "johnbrown@example.net", "+1-375-596-4328x74145")
    ,@("Patty", "Gregory", "sheltonphillip@example.com", "496-302-7768x50249")
    ,@("Albert", "Smith", "charles87@example.org", "+1-230-864-2679x949")
    ,@("Jillian", "Lang", "garrettkelsey@example.net", "419.981.6103x4274")
    ,@("Nancy", "Perez", "sbrown@example.net", "+1-568-329-1081x978")
    ,@("Lisa", "Benson", "danduncan@example.org", "(934)282-5500")
    ,@("Caleb", "Torres", "rmccarthy@example.org", "755.892.7111x981")
    ,@("Jenna", "Gonzales", "ajohnson@example.com", "(254)512-0559x5347")
    ,@("Larry", "Mclaughlin", "suttondavid@example.net", "(994)618-9510x19469")
    ,@("Shelby", "Miller", "chad11@example.net", "(782)413-8582")
    ,@("Pamela", "Castillo", "ohoffman@example.org", "+1-818-595-2652x457")
    ,@("William", "Thomas", "egraves@example.org", "270-272-3681")
    ,@("Lauren", "Jimenez", "hutchinsonnicole@example.org", "001-212-839-8551")
    ,@("Justin", "Jenkins", "ellenferguson@example.net", "764.415.3410x62196")
    ,@("Christopher", "Long", "schmidtmichelle@example.org", "674.921.0391x576")
    ,@("Jose", "Harrington", "wmercer@example.net", "(740)735-3357x46828")
    ,@("Erica", "Adams", "carrie65@example.org", "429.791.1693")
    ,@("Teresa", "Perry", "davidwilson@example.org", "217-298-8512x66487")
)

$ws.Range("A1:D49").NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 1
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
